$d = $word.ActiveDocument

# Step 1: bookmark at the end of paragraph 3 (after the hyperlink, before paragraph mark)
$r = $d.Content
$r.Find.Execute("valuevsreftypes.aspx") | Out-Null
$r.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r) | Out-Null

# Step 2: insert new paragraph with hyperlink after paragraph 3, using InsertXML for exact control
$p3 = $d.Paragraphs.Item(3)
$sel = $word.Selection
$sel.SetRange($p3.Range.End, $p3.Range.End)
$sel.Collapse(0)
$url = "https://www.tutlane.com/tutorial/csharp/csharp-value-type-and-reference-type-with-examples"
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:pPr><w:spacing w:after="160"/></w:pPr><w:hyperlink r:id="rIdNEWHL" w:history="1"><w:r><w:rPr><w:color w:val="0000FF"/><w:u w:val="single"/></w:rPr><w:t>' + $url + '</w:t></w:r></w:hyperlink></w:p><w:p><w:r><w:t>PLACEHOLDER</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/_rels/document.xml.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rIdNEWHL" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="' + $url + '" TargetMode="External"/></Relationships></pkg:xmlData></pkg:part></pkg:package>'
$sel.InsertXML($xml)
Write-Host "done"
